# testDataSource.xlsx edit: add a "pADATremap"/"nextAbbrevOverlay" lookup table,
# drop Sheet3, tidy up Sheet1/Sheet2 formatting and register a LookupTable name.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove Sheet3 -----------------------------------------------------
$wb.Worksheets.Item("Sheet3").Delete()

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 cell values --------------------------------------------------
$ws1.Range("A1").Value = "basePath"
$ws1.Range("B1").Value = "y:\adat"

$ws1.Range("A2").Value = "MRIcroNexe"
$ws1.Range("B2").Value = "C:\Program Files\mricron\MRIcroN.exe"

$ws1.Range("A3").Value = "pADATremap"
$ws1.Range("B3").Value = ".\masked_roi98_mniwholebrain_fromspm_wroi99_wholecube_both_p-overlay_adathreshold_remap_clustercorrected.hdr"
$ws1.Range("C3").Value = "red"

$ws1.Range("A4").Value = "nextAbbrevOverlay"
$ws1.Range("B4").Value = ".\next\goes\here"
$ws1.Range("C4").Value = "green"

# old row 5 ("again" / "something") is gone -- clear it out
$ws1.Range("A5:C5").Clear()

# --- Sheet1 formatting ----------------------------------------------------
# Column A: bold white-on-accent text, right aligned (key column header look)
$colA = $ws1.Range("A1:A4")
$colA.Interior.ThemeColor = 10
$colA.Font.ThemeColor = 2
$colA.Font.Bold = $true
$colA.HorizontalAlignment = -4152

# Column B: left aligned with an indent (value column)
$colB = $ws1.Range("B1:B4")
$colB.Font.Bold = $false
$colB.HorizontalAlignment = -4131
$colB.IndentLevel = 1

# Column widths
$ws1.Columns("A").ColumnWidth = 17.08984375
$ws1.Columns("B").ColumnWidth = 106.6328125

$ws1.PageSetup.Orientation = 1

# --- Sheet2 ---------------------------------------------------------------
$ws2.Columns("A").ColumnWidth = 86.453125

# --- Defined name ----------------------------------------------------------
$wb.Names.Add("LookupTable", "=Sheet1!`$A`$1:`$C`$4")

# make sure Sheet1 stays the active sheet/tab, with the same selection Excel left it in
$ws1.Activate()
$ws1.Range("B10").Select()
